# Fruta / hortaliza, semanal
# Insert two new data rows (157 and 158) into the "Frutilla" sheet, shifting
# all existing rows 157:263 down to 159:265. The two newly inserted rows get
# brand new price-report data; everything else keeps its original content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 157 - this shifts rows
# 157:263 down to 159:265 and extends the sheet dimension accordingly.
$ws.Rows("157:158").Insert()

# ---- New row 157 ----
$ws.Range("A157").Value = 7
$ws.Range("B157").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C157").Value = "Ñuble"
$ws.Range("D157").Value = 44606
$ws.Range("E157").Value = 16
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100101
$ws.Range("H157").Value = "Berries"
$ws.Range("I157").Value = 100112025
$ws.Range("J157").Value = "Frutilla"
$ws.Range("K157").Value = "Sin especificar"
$ws.Range("L157").Value = "Especial"
$ws.Range("M157").Value = 80
$ws.Range("N157").Value = 7000
$ws.Range("O157").Value = 7000
$ws.Range("P157").Value = 7000
$ws.Range("Q157").Value = "$/bandeja 7 kilos"
$ws.Range("R157").Value = "Provincia de Diguillín"
$ws.Range("S157").Value = 1000
$ws.Range("T157").Value = 7

# ---- New row 158 ----
$ws.Range("A158").Value = 7
$ws.Range("B158").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C158").Value = "Ñuble"
$ws.Range("D158").Value = 44606
$ws.Range("E158").Value = 16
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100101
$ws.Range("H158").Value = "Berries"
$ws.Range("I158").Value = 100112025
$ws.Range("J158").Value = "Frutilla"
$ws.Range("K158").Value = "Sin especificar"
$ws.Range("L158").Value = "Primera"
$ws.Range("M158").Value = 120
$ws.Range("N158").Value = 6000
$ws.Range("O158").Value = 6500
$ws.Range("P158").Value = 6250
$ws.Range("Q158").Value = "$/bandeja 7 kilos"
$ws.Range("R158").Value = "Provincia de Diguillín"
$ws.Range("S158").Value = 893
$ws.Range("T158").Value = 7
